$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Rows.Item(82).Insert()

$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K137"))

$ws.Range("A81:K81").Copy()
$ws.Range("A82:K82").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("G82").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G137").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

$ws.Range("B82").Value = "SL(1-0-0)"
$ws.Range("H82").Value = 1
$ws.Range("K82").Value = [DateTime]::FromOADate(44910)

$ws.Range("B83").Select()
